$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 258, shifting the
# existing rows 258-354 down to 260-356 (same layout/style carried along).
$ws.Rows("258:259").Insert()

# Shared metadata values that are constant across every data row in this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108006
$categoria = "Plátano"
$variedad  = "Sin especificar"
$unidad    = "$/caja 20 kilos"
$origen    = "Ecuador"
$kgUnidad  = 20

# New row 258: Pintón, 2021-09-29
$r = 258
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44468
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Pintón"
$ws.Cells.Item($r,13).Value = 180
$ws.Cells.Item($r,14).Value = 13000
$ws.Cells.Item($r,15).Value = 13000
$ws.Cells.Item($r,16).Value = 13000
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 650
$ws.Cells.Item($r,20).Value = $kgUnidad

# New row 259: Primera Pintón, 2021-09-29
$r = 259
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44468
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Primera Pintón"
$ws.Cells.Item($r,13).Value = 500
$ws.Cells.Item($r,14).Value = 14000
$ws.Cells.Item($r,15).Value = 15000
$ws.Cells.Item($r,16).Value = 14500
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 725
$ws.Cells.Item($r,20).Value = $kgUnidad
